$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 34.71251733333333
$ws.Range("H2").Value = 104.137552
$ws.Range("I2").Value = 0.111750244749681
$ws.Range("J2").Value = 0.111750244749681
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1528053333333333
$ws.Range("N2").Value = 0.458416
$ws.Range("O2").Value = 0.01103433215988526
$ws.Range("P2").Value = 0.01103433215988526
$ws.Range("Q2").Value = 5.304257781959111
$ws.Range("R2").Value = 47.738320037632
$ws.Range("S2").Value = 0.001233089319516454
$ws.Range("T2").Value = 0.001233089319516454
$ws.Range("G3").Value = 34.71251733333333
$ws.Range("H3").Value = 104.137552
$ws.Range("I3").Value = 0.111750244749681
$ws.Range("J3").Value = 0.111750244749681
$ws.Range("O3").Value = 0.8539197603380489
$ws.Range("P3").Value = 0.8539197603380488
$ws.Range("Q3").Value = 410.4834319206187
$ws.Range("R3").Value = 3694.350887285568
$ws.Range("S3").Value = 0.09542574221436594
$ws.Range("T3").Value = 0.09542574221436592
$ws.Range("G4").Value = 34.71251733333333
$ws.Range("H4").Value = 104.137552
$ws.Range("I4").Value = 0.111750244749681
$ws.Range("J4").Value = 0.111750244749681
$ws.Range("M4").Value = 1.712948333333333
$ws.Range("N4").Value = 5.138845
$ws.Range("O4").Value = 0.1236949029880405
$ws.Range("P4").Value = 0.1236949029880405
$ws.Range("Q4").Value = 59.46074871193778
$ws.Range("R4").Value = 535.14673840744
$ws.Range("S4").Value = 0.01382293568320158
$ws.Range("T4").Value = 0.01382293568320157
$ws.Range("G5").Value = 34.71251733333333
$ws.Range("H5").Value = 104.137552
$ws.Range("I5").Value = 0.111750244749681
$ws.Range("J5").Value = 0.111750244749681
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1571906666666667
$ws.Range("N5").Value = 0.471572
$ws.Range("O5").Value = 0.01135100451402528
$ws.Range("P5").Value = 0.01135100451402528
$ws.Range("Q5").Value = 5.456483741304889
$ws.Range("R5").Value = 49.108353671744
$ws.Range("S5").Value = 0.001268477532597059
$ws.Range("T5").Value = 0.001268477532597059
$ws.Range("H6").Value = 578.4917909999999
$ws.Range("I6").Value = 0.620780861354714
$ws.Range("J6").Value = 0.6207808613547139
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1528053333333333
$ws.Range("N6").Value = 0.458416
$ws.Range("O6").Value = 0.01103433215988526
$ws.Range("P6").Value = 0.01103433215988526
$ws.Range("Q6").Value = 29.46554365145066
$ws.Range("R6").Value = 265.1898928630559
$ws.Range("S6").Value = 0.006849902222687594
$ws.Range("T6").Value = 0.006849902222687592
$ws.Range("H7").Value = 578.4917909999999
$ws.Range("I7").Value = 0.620780861354714
$ws.Range("J7").Value = 0.6207808613547139
$ws.Range("O7").Value = 0.8539197603380489
$ws.Range("P7").Value = 0.8539197603380488
$ws.Range("R7").Value = 20522.39197411004
$ws.Range("S7").Value = 0.530097044350465
$ws.Range("T7").Value = 0.5300970443504648
$ws.Range("H8").Value = 578.4917909999999
$ws.Range("I8").Value = 0.620780861354714
$ws.Range("J8").Value = 0.6207808613547139
$ws.Range("M8").Value = 1.712948333333333
$ws.Range("N8").Value = 5.138845
$ws.Range("O8").Value = 0.1236949029880405
$ws.Range("P8").Value = 0.1236949029880405
$ws.Range("Q8").Value = 330.3088497468216
$ws.Range("R8").Value = 2972.779647721395
$ws.Range("S8").Value = 0.07678742842210357
$ws.Range("T8").Value = 0.07678742842210355
$ws.Range("H9").Value = 578.4917909999999
$ws.Range("I9").Value = 0.620780861354714
$ws.Range("J9").Value = 0.6207808613547139
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1571906666666667
$ws.Range("N9").Value = 0.471572
$ws.Range("O9").Value = 0.01135100451402528
$ws.Range("P9").Value = 0.01135100451402528
$ws.Range("Q9").Value = 30.31117009616133
$ws.Range("R9").Value = 272.8005308654519
$ws.Range("S9").Value = 0.00704648635945786
$ws.Range("T9").Value = 0.007046486359457858
$ws.Range("G10").Value = 19.96051866666667
$ws.Range("H10").Value = 59.881556
$ws.Range("I10").Value = 0.06425903442584988
$ws.Range("J10").Value = 0.06425903442584986
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1528053333333333
$ws.Range("N10").Value = 0.458416
$ws.Range("O10").Value = 0.01103433215988526
$ws.Range("P10").Value = 0.01103433215988526
$ws.Range("Q10").Value = 3.050073708366222
$ws.Range("R10").Value = 27.450663375296
$ws.Range("S10").Value = 0.0007090555301283294
$ws.Range("T10").Value = 0.0007090555301283292
$ws.Range("G11").Value = 19.96051866666667
$ws.Range("H11").Value = 59.881556
$ws.Range("I11").Value = 0.06425903442584988
$ws.Range("J11").Value = 0.06425903442584986
$ws.Range("O11").Value = 0.8539197603380489
$ws.Range("P11").Value = 0.8539197603380488
$ws.Range("Q11").Value = 236.0376842315894
$ws.Range("R11").Value = 2124.339158084304
$ws.Range("S11").Value = 0.05487205927647616
$ws.Range("T11").Value = 0.05487205927647614
$ws.Range("G12").Value = 19.96051866666667
$ws.Range("H12").Value = 59.881556
$ws.Range("I12").Value = 0.06425903442584988
$ws.Range("J12").Value = 0.06425903442584986
$ws.Range("M12").Value = 1.712948333333333
$ws.Range("N12").Value = 5.138845
$ws.Range("O12").Value = 0.1236949029880405
$ws.Range("P12").Value = 0.1236949029880405
$ws.Range("Q12").Value = 34.19133718253556
$ws.Range("R12").Value = 307.72203464282
$ws.Range("S12").Value = 0.007948515029410656
$ws.Range("T12").Value = 0.007948515029410652
$ws.Range("G13").Value = 19.96051866666667
$ws.Range("H13").Value = 59.881556
$ws.Range("I13").Value = 0.06425903442584988
$ws.Range("J13").Value = 0.06425903442584986
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1571906666666667
$ws.Range("N13").Value = 0.471572
$ws.Range("O13").Value = 0.01135100451402528
$ws.Range("P13").Value = 0.01135100451402528
$ws.Range("Q13").Value = 3.137607236225778
$ws.Range("R13").Value = 28.238465126032
$ws.Range("S13").Value = 0.0007294045898347278
$ws.Range("T13").Value = 0.0007294045898347275
$ws.Range("G14").Value = 63.12224
$ws.Range("H14").Value = 189.36672
$ws.Range("I14").Value = 0.2032098594697551
$ws.Range("J14").Value = 0.2032098594697551
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1528053333333333
$ws.Range("N14").Value = 0.458416
$ws.Range("O14").Value = 0.01103433215988526
$ws.Range("P14").Value = 0.01103433215988526
$ws.Range("Q14").Value = 9.645414923946666
$ws.Range("R14").Value = 86.80873431552
$ws.Range("S14").Value = 0.002242285087552883
$ws.Range("T14").Value = 0.002242285087552883
$ws.Range("G15").Value = 63.12224
$ws.Range("H15").Value = 189.36672
$ws.Range("I15").Value = 0.2032098594697551
$ws.Range("J15").Value = 0.2032098594697551
$ws.Range("O15").Value = 0.8539197603380489
$ws.Range("P15").Value = 0.8539197603380488
$ws.Range("Q15").Value = 746.43487987072
$ws.Range("R15").Value = 6717.91391883648
$ws.Range("S15").Value = 0.1735249144967419
$ws.Range("T15").Value = 0.1735249144967419
$ws.Range("G16").Value = 63.12224
$ws.Range("H16").Value = 189.36672
$ws.Range("I16").Value = 0.2032098594697551
$ws.Range("J16").Value = 0.2032098594697551
$ws.Range("M16").Value = 1.712948333333333
$ws.Range("N16").Value = 5.138845
$ws.Range("O16").Value = 0.1236949029880405
$ws.Range("P16").Value = 0.1236949029880405
$ws.Range("Q16").Value = 108.1251358042667
$ws.Range("R16").Value = 973.1262222383999
$ws.Range("S16").Value = 0.02513602385332471
$ws.Range("T16").Value = 0.0251360238533247
$ws.Range("G17").Value = 63.12224
$ws.Range("H17").Value = 189.36672
$ws.Range("I17").Value = 0.2032098594697551
$ws.Range("J17").Value = 0.2032098594697551
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.1571906666666667
$ws.Range("N17").Value = 0.471572
$ws.Range("O17").Value = 0.01135100451402528
$ws.Range("P17").Value = 0.01135100451402528
$ws.Range("Q17").Value = 9.922226987093334
$ws.Range("R17").Value = 89.30004288383999
$ws.Range("S17").Value = 0.002306636032135633
$ws.Range("T17").Value = 0.002306636032135633
